$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.565.21"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.41%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.652.80"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.02%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.31"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.91%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.88"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.52%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.631"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +7.29%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.43%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.396"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.62%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.37%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.51%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.84"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.15%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.92%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.129.41"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.06%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.426.28"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.32%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.645.92"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.16%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.63"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.07%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.77"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.29%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.59%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "349.14"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.01%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.00%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.99"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.87%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.42%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.62"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.40%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.63%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.80%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.163"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.25%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.96"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.67%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "540.42"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.11%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.34%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.11"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.08%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.57%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.42"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.27%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.42"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.16%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.52%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.36"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.61%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.01%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "156.06"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.76%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.89%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.02%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "161.18"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.15%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.71%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.30%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0605"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.59%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.55"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.52%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.24%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.35%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.100"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.45%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆0252"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +6.52%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.31%  "
